# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Sun Apr 30 10:01:07 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.574.34"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.925.58"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.48"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4824"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4059"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08217"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.82"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.920.81"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.115"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.313"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.68"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06866"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.69"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.579.68"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.676"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.04"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.178"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.138.56"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.17"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.426"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.097"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.71"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.013"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09600"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.613"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.563"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.382"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06373"
$ws.Range("E36").Value = "  +4.47%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5964"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.74"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.890"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1849"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.503"
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.285"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07544"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5564"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.979"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "119.13"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.437"
$ws.Range("E51").Value = "  +0.92%  "
